$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric must be forced to remain text
# (matching the original inline-string cell type) by using a leading
# apostrophe, then the style is reset to Normal so no stray number-format
# style gets attached to the cell.

$ws.Range("D2").Value = "'243.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.418"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05936"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.453"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.534"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = "'0.9124"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.07478"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03300"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03059"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09352"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.850"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001591"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04666"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006076"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.004990"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.0009844"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001100"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.606"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.136"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01123"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23OneONEBestin24h"
$ws.Range("D25").Value = "'0.3228"
$ws.Range("D25").Style = "Normal"
$ws.Range("D40").Value = "'0.03947"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006210"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Value = "'0.008682"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005231"
$ws.Range("D45").Style = "Normal"
